$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "predicted results" table (rows 11-13, column C = Edge)
$ws.Range("C11").Value = 0.50314499999999995
$ws.Range("C12").Value = 0.42375699999999999
$ws.Range("C13").Value = 0.31222

# Update the view: scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("B19").Select()
